# The sheet currently stores each product "record" smeared across 5 rows
# (one populated column - B, C, D, E, F respectively - per row). Consolidate
# each group of 5 rows into a single row with all five fields populated,
# re-number the leading index column (A) as 0..4, and drop the now-empty
# trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$groupCount = 5
$colLetters = @("B", "C", "D", "E", "F")

# First, read every value out of its current (scattered) position before any
# writes/deletes happen, so row addressing stays stable while we read.
$values = @{}
for ($g = 0; $g -lt $groupCount; $g++) {
    $baseRow = 2 + ($g * 5)
    for ($k = 0; $k -lt 5; $k++) {
        $col = $colLetters[$k]
        $srcRow = $baseRow + $k
        $addr = "$col$srcRow"
        $values[[string]($g.ToString() + "_" + $col)] = $ws.Range($addr).Value()
    }
}

# Now write the consolidated rows 2..6.
for ($g = 0; $g -lt $groupCount; $g++) {
    $targetRow = 2 + $g
    $ws.Range("A$targetRow").Value = $g
    foreach ($col in $colLetters) {
        $val = $values[[string]($g.ToString() + "_" + $col)]
        $ws.Range("$col$targetRow").Value = $val
    }
}

# Finally, remove the now-redundant trailing rows (previously rows 7..26,
# which held the scattered data we already folded into rows 2..6).
$ws.Range("A7:A26").EntireRow.Delete()
